$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.683.04"
$ws.Range("E2").Value = "  -1.48%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.596.72"
$ws.Range("E3").Value = "  -1.65%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.04"
$ws.Range("E5").Value = "  -1.31%  "

# Row 6
$ws.Range("E6").Value = "  -0.62%  "

# Row 7
$ws.Range("E7").Value = "  +0.12%  "

# Row 8
$ws.Range("E8").Value = "  -1.53%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0618"
$ws.Range("E9").Value = "  -2.37%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.74"
$ws.Range("E10").Value = "  -1.94%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0837"
$ws.Range("E11").Value = "  -1.25%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.819.95"
$ws.Range("E12").Value = "  -1.67%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.594.07"
$ws.Range("E13").Value = "  -1.50%  "

# Row 14
$ws.Range("E14").Value = "  -2.86%  "

# Row 15
$ws.Range("E15").Value = "  -3.27%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.12"
$ws.Range("E16").Value = "  +0.39%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.646.01"
$ws.Range("E17").Value = "  -1.42%  "

# Row 18
$ws.Range("E18").Value = "  -1.89%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "210.31"
$ws.Range("E19").Value = "  -1.83%  "

# Row 20
$ws.Range("E20").Value = "  +0.16%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.70"
$ws.Range("E21").Value = "  -2.23%  "

# Row 22
$ws.Range("E22").Value = "  -2.65%  "

# Row 23
$ws.Range("E23").Value = "  -2.20%  "

# Row 24
$ws.Range("E24").Value = "  -1.89%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.62"

# Row 26
$ws.Range("E26").Value = "  +0.00%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.20"
$ws.Range("E27").Value = "  -2.64%  "

# Row 28
$ws.Range("E28").Value = "  -0.26%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.34"
$ws.Range("E29").Value = "  -1.78%  "

# Row 30
$ws.Range("E30").Value = "  -1.64%  "

# Row 31
$ws.Range("E31").Value = "  -1.44%  "

# Row 32
$ws.Range("E32").Value = "  -3.53%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.675"
$ws.Range("E33").Value = "  -10.07%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.93"
$ws.Range("E34").Value = "  -2.82%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.291.29"
$ws.Range("E35").Value = "  -4.81%  "

# Row 36
$ws.Range("E36").Value = "  -0.77%  "

# Row 37
$ws.Range("E37").Value = "  -5.75%  "

# Row 38
$ws.Range("E38").Value = "  -3.45%  "

# Row 39
$ws.Range("E39").Value = "  -1.45%  "

# Row 40
$ws.Range("E40").Value = "  +0.10%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.40"
$ws.Range("E41").Value = "  +0.79%  "

# Row 42
$ws.Range("E42").Value = "  -1.62%  "

# Row 43
$ws.Range("E43").Value = "  -1.36%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.81"
$ws.Range("E44").Value = "  -2.07%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.732.13"
$ws.Range("E45").Value = "  -1.73%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.85"
$ws.Range("E46").Value = "  -0.38%  "

# Row 47
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.872"
$ws.Range("E47").Value = "  -1.46%  "

# Row 48
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.63"
$ws.Range("E48").Value = "  -1.35%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0986"
$ws.Range("E49").Value = "  -3.40%  "

# Row 50
$ws.Range("E50").Value = "  -2.14%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.49"
$ws.Range("E51").Value = "  -2.39%  "
